# Fix Training Data Issue (#48)
# The "Date" column (BF) held the literal folder/file label "6-11-2007-08"
# instead of an actual date string. Because NBA.com stats pages are
# labelled by the day *after* the games they report, the data was
# effectively tagged one day off. Correct every data row's Date cell
# (BF2:BF31) to read "2008-06-11" as plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "6-11-2007-08"
$newDate = "2008-06-11"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Text -eq $oldDate) {
        # Force text so Excel doesn't auto-convert the ISO-looking
        # "2008-06-11" into a date serial number, then drop the
        # temporary number-format override so the cell keeps its
        # original (default) style.
        $cell.NumberFormat = "@"
        $cell.Value = $newDate
        $cell.ClearFormats()
    }
}
